# Updates cryptos list values (Price column D, Volume(1h) column E)
# per the commit "Updated cryptos list on Sun Sep 24 22:50:50 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells: force text storage via a temporary Text number format
# so Excel does not silently convert these string-formatted prices into doubles
# (which would lose trailing zeros / exact formatting), then clear the format again
# so the cell keeps its original (default) style, matching the source file.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '26.617.50'
$ws.Range('E2').Value = '  -0.37%  '
Set-TextValue 'D3' '1.595.74'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue 'D5' '210.73'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('E10').Value = '  +0.29%  '
Set-TextValue 'D11' '0.0844'
$ws.Range('E11').Value = '  +0.18%  '
Set-TextValue 'D12' '1.819.49'
Set-TextValue 'D13' '1.592.56'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('E15').Value = '  -0.28%  '
Set-TextValue 'D16' '64.50'
$ws.Range('E16').Value = '  -1.23%  '
Set-TextValue 'D17' '26.592.47'
$ws.Range('E17').Value = '  -0.36%  '
Set-TextValue 'D18' '0.0₃0738'
$ws.Range('E18').Value = '  -2.01%  '
$ws.Range('E19').Value = '  -0.06%  '
Set-TextValue 'D20' '208.64'
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('E21').Value = '  -2.34%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('E24').Value = '  +0.20%  '
Set-TextValue 'D25' '144.90'
$ws.Range('E25').Value = '  +1.85%  '
Set-TextValue 'D26' '1.01'
$ws.Range('E26').Value = '  +0.02%  '
Set-TextValue 'D29' '15.29'
$ws.Range('E29').Value = '  -0.34%  '
Set-TextValue 'D30' '0.0507'
$ws.Range('E30').Value = '  -2.72%  '
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('E33').Value = '  -0.41%  '
Set-TextValue 'D34' '1.282.65'
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').Value = '  +0.34%  '
Set-TextValue 'D36' '1.23'
$ws.Range('E36').Value = '  +11.94%  '
Set-TextValue 'D37' '0.601'
$ws.Range('E37').Value = '  -3.56%  '
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('E42').Value = '  -1.81%  '
Set-TextValue 'D43' '0.769'
$ws.Range('E43').Value = '  -1.85%  '
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('E45').Value = '  -0.32%  '
Set-TextValue 'D46' '89.31'
$ws.Range('E46').Value = '  -2.06%  '
$ws.Range('E48').Value = '  +2.49%  '
$ws.Range('E49').Value = '  +0.54%  '
Set-TextValue 'D50' '7.51'
$ws.Range('E50').Value = '  +1.86%  '
$ws.Range('E51').Value = '  -0.04%  '
